# Replicate the "Fondling" style block but for Offense = "Criminal Offenses -
# Non Forcible Sex Offenses", with data only for Date = "sum2013".
# This appends 36 new rows (722-757) to Sheet1: 4 reporting-location blocks
# (On Campus excl. Res Halls / On Campus Res Halls / Non-Campus / Public
# Property) x 9 institution sectors each.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$offense = "Criminal Offenses - Non Forcible Sex Offenses"
$dateLabel = "sum2013"

$sectors = @(
  "Public, 4-year or above",
  "Private nonprofit, 4-year or above",
  "Private for-profit, 4-year or above",
  "Public, 2-year",
  "Private nonprofit, 2-year",
  "Private for-profit, 2-year",
  "Public, less-than 2-year",
  "Private nonprofit, less-than 2-year",
  "Private for-profit, less-than 2-year"
)

$locations = @(
  "On Campus (excluding Residence Halls)",
  "On Campus (Residence Halls)",
  "Non-Campus",
  "Public Property"
)

# Counts (Column E), in row order 722..757, 9 per location block.
$counts = @(
  9, 8, 2, 7, 0, 0, 1, 0, 0,
  9, 4, 0, 4, 2, 0, 0, 0, 0,
  1, 1, 2, 5, 0, 0, 0, 0, 0,
  0, 1, 0, 1, 0, 0, 0, 0, 0
)

$startRow = 722
$row = $startRow

for ($locIdx = 0; $locIdx -lt $locations.Count; $locIdx++) {
  $loc = $locations[$locIdx]
  for ($secIdx = 0; $secIdx -lt $sectors.Count; $secIdx++) {
    $sec = $sectors[$secIdx]
    $countIdx = $row - $startRow
    $val = $counts[$countIdx]

    $ws.Cells.Item($row, 1).Value2 = $sec
    $ws.Cells.Item($row, 2).Value2 = $loc
    $ws.Cells.Item($row, 3).Value2 = $offense
    $ws.Cells.Item($row, 4).Value2 = $dateLabel
    $ws.Cells.Item($row, 5).Value2 = $val

    $row = $row + 1
  }
}

$endRow = $row - 1

# Apply the same cell formatting (style "1" with quotePrefix, as used by the
# rest of the table) that row 721 already has, across the newly added block.
$ws.Range("A721:E721").Copy()
$ws.Range("A" + $startRow + ":E" + $endRow).PasteSpecial(-4122)

# Replicate the pre-existing stray styled-but-empty cell in column F that
# appears at the first "On Campus (Residence Halls)" / sum2013 row of every
# offense block (see e.g. F137, F245, F353, F461, F551, F596, F668).
$ws.Range("F668").Copy()
$ws.Range("F731").PasteSpecial(-4122)

# Update the selection to match the end-state of the edit.
$ws.Range("G751").Select()
